$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before the existing "Date" column (E), pushing it to F.
$ws.Range("E1").EntireColumn.Insert()

# 2. Populate the new column E: header "Data" plus "Historical"/"Forecast" tags per row.
$ws.Range("E1").Value = "Data"

for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 5).Value = "Historical"
}

for ($row = 48; $row -le 71; $row++) {
    $ws.Cells.Item($row, 5).Value = "Forecast"
}

# 3. Update the revised Sale values (column C) for the 2016 rows.
$newSaleValues = @{
    48 = 2761958
    49 = 2656665.25
    50 = 3099057.75
    51 = 2873607.25
    52 = 3327835.25
    53 = 3356062
    54 = 3391942.75
    55 = 2991382.5
    56 = 2664295.25
    57 = 2588209.75
    58 = 2702838.25
    59 = 2761943.25
    60 = 21539936.0074994
    61 = 20413770.6013595
    62 = 24325953.0976278
    63 = 22993466.3485849
    64 = 26691951.4191559
    65 = 26989964.0105518
    66 = 26948630.7647638
    67 = 24091579.3491059
    68 = 20523492.4086428
    69 = 20011748.6685998
    70 = 21177435.4858385
    71 = 20855799.1096099
}

foreach ($row in $newSaleValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newSaleValues[$row]
}
